$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "18×21=" -> "46×67="
$cell = $t.Cell(1, 1)
$cell.Range.Text = "46×67="

# Row 1, Col 2: "34×82=" -> "69×46="
$cell = $t.Cell(1, 2)
$cell.Range.Text = "69×46="

# Row 1, Col 3: "49×70=" -> "20×94="
$cell = $t.Cell(1, 3)
$cell.Range.Text = "20×94="

# Row 1, Col 4: "62×16=" -> "37×13="
$cell = $t.Cell(1, 4)
$cell.Range.Text = "37×13="

# Row 1, Col 5: "45×36=" -> "29×60="
$cell = $t.Cell(1, 5)
$cell.Range.Text = "29×60="

# Row 5, Col 1: "87×98=" -> "86×99="
$cell = $t.Cell(5, 1)
$cell.Range.Text = "86×99="

# Row 5, Col 2: "63×84=" -> "92×57="
$cell = $t.Cell(5, 2)
$cell.Range.Text = "92×57="

# Row 5, Col 3: "95×18=" -> "16×32="
$cell = $t.Cell(5, 3)
$cell.Range.Text = "16×32="

# Row 5, Col 4: "40×70=" -> "46×52="
$cell = $t.Cell(5, 4)
$cell.Range.Text = "46×52="

# Row 5, Col 5: "32×67=" -> "55×21="
$cell = $t.Cell(5, 5)
$cell.Range.Text = "55×21="

# Row 10, Col 1: "29×68=" -> "64×81="
$cell = $t.Cell(10, 1)
$cell.Range.Text = "64×81="

# Row 10, Col 2: "16×63=" -> "21×68="
$cell = $t.Cell(10, 2)
$cell.Range.Text = "21×68="

# Row 10, Col 3: "37×75=" -> "34×23="
$cell = $t.Cell(10, 3)
$cell.Range.Text = "34×23="

# Row 10, Col 4: "49×19=" -> "60×58="
$cell = $t.Cell(10, 4)
$cell.Range.Text = "60×58="

# Row 10, Col 5: "54×60=" -> "68×56="
$cell = $t.Cell(10, 5)
$cell.Range.Text = "68×56="

# Row 15, Col 1: "36×25=" -> "72×32="
$cell = $t.Cell(15, 1)
$cell.Range.Text = "72×32="

# Row 15, Col 2: "48×44=" -> "49×29="
$cell = $t.Cell(15, 2)
$cell.Range.Text = "49×29="

# Row 15, Col 3: "49×70=" -> "28×36="
$cell = $t.Cell(15, 3)
$cell.Range.Text = "28×36="

# Row 15, Col 4: "55×67=" -> "59×31="
$cell = $t.Cell(15, 4)
$cell.Range.Text = "59×31="

# Row 15, Col 5: "26×87=" -> "68×26="
$cell = $t.Cell(15, 5)
$cell.Range.Text = "68×26="

# Row 20, Col 1: "52×75=" -> "95×59="
$cell = $t.Cell(20, 1)
$cell.Range.Text = "95×59="

# Row 20, Col 2: "83×67=" -> "51×11="
$cell = $t.Cell(20, 2)
$cell.Range.Text = "51×11="

# Row 20, Col 3: "92×76=" -> "46×97="
$cell = $t.Cell(20, 3)
$cell.Range.Text = "46×97="

# Row 20, Col 4: "53×71=" -> "65×36="
$cell = $t.Cell(20, 4)
$cell.Range.Text = "65×36="

# Row 20, Col 5: "27×54=" -> "65×85="
$cell = $t.Cell(20, 5)
$cell.Range.Text = "65×85="
